{"js": "const table = context.document.body.tables.getFirst();\n\n// Map of (row, col) -> new text, applied by absolute table position so\n// that the replacements are unambiguous even when a new value coincides\n// with another cell's original value.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"95\u00f79=\", newText: \"66\u00f75=\" },\n  { row: 0, col: 1, oldText: \"56\u00f73=\", newText: \"62\u00f73=\" },\n  { row: 0, col: 2, oldText: \"43\u00f79=\", newText: \"57\u00f73=\" },\n  { row: 0, col: 3, oldText: \"74\u00f78=\", newText: \"87\u00f73=\" },\n  { row: 0, col: 4, oldText: \"41\u00f78=\", newText: \"23\u00f76=\" },\n  { row: 4, col: 0, oldText: \"50\u00f75=\", newText: \"66\u00f73=\" },\n  { row: 4, col: 1, oldText: \"84\u00f73=\", newText: \"50\u00f74=\" },\n  { row: 4, col: 2, oldText: \"25\u00f77=\", newText: \"43\u00f73=\" },\n  { row: 4, col: 3, oldText: \"87\u00f74=\", newText: \"43\u00f79=\" },\n  { row: 4, col: 4, oldText: \"36\u00f75=\", newText: \"33\u00f73=\" },\n  { row: 8, col: 0, oldText: \"28\u00f74=\", newText: \"15\u00f75=\" },\n  { row: 8, col: 1, oldText: \"82\u00f76=\", newText: \"99\u00f72=\" },\n  { row: 8, col: 2, oldText: \"70\u00f74=\", newText: \"78\u00f72=\" },\n  { row: 8, col: 3, oldText: \"62\u00f73=\", newText: \"66\u00f78=\" },\n  { row: 8, col: 4, oldText: \"36\u00f77=\", newText: \"88\u00f74=\" },\n  { row: 12, col: 0, oldText: \"57\u00f75=\", newText: \"96\u00f72=\" },\n  { row: 12, col: 1, oldText: \"88\u00f74=\", newText: \"28\u00f78=\" },\n  { row: 12, col: 2, oldText: \"74\u00f72=\", newText: \"73\u00f77=\" },\n  { row: 12, col: 3, oldText: \"68\u00f77=\", newText: \"81\u00f73=\" },\n  { row: 12, col: 4, oldText: \"90\u00f72=\", newText: \"75\u00f73=\" },\n  { row: 16, col: 0, oldText: \"26\u00f76=\", newText: \"33\u00f77=\" },\n  { row: 16, col: 1, oldText: \"28\u00f76=\", newText: \"54\u00f79=\" },\n  { row: 16, col: 2, oldText: \"96\u00f73=\", newText: \"69\u00f76=\" },\n  { row: 16, col: 3, oldText: \"44\u00f79=\", newText: \"47\u00f76=\" },\n  { row: 16, col: 4, oldText: \"41\u00f74=\", newText: \"52\u00f76=\" },\n];\n\n// Load current values first so we can confirm each cell still holds the\n// expected \"before\" text before overwriting it (guards against the table\n// layout having shifted for some unrelated reason).\nconst cells = replacements.map(({ row, col }) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nreplacements.forEach(({ oldText, newText }, i) => {\n  const cell = cells[i];\n  if (cell.value !== oldText) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: expected \"${oldText}\" but found \"${cell.value}\"`\n    );\n  }\n  cell.value = newText;\n});\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Map of (row, col) -> new text, applied by absolute table position (1-based,\n# matching Word COM's Table.Cell(row, col)) so the replacements stay unambiguous\n# even when a new value coincides with another cell's original value.\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"95\u00f79=\"; NewText = \"66\u00f75=\" }\n    @{ Row = 1; Col = 2; OldText = \"56\u00f73=\"; NewText = \"62\u00f73=\" }\n    @{ Row = 1; Col = 3; OldText = \"43\u00f79=\"; NewText = \"57\u00f73=\" }\n    @{ Row = 1; Col = 4; OldText = \"74\u00f78=\"; NewText = \"87\u00f73=\" }\n    @{ Row = 1; Col = 5; OldText = \"41\u00f78=\"; NewText = \"23\u00f76=\" }\n    @{ Row = 5; Col = 1; OldText = \"50\u00f75=\"; NewText = \"66\u00f73=\" }\n    @{ Row = 5; Col = 2; OldText = \"84\u00f73=\"; NewText = \"50\u00f74=\" }\n    @{ Row = 5; Col = 3; OldText = \"25\u00f77=\"; NewText = \"43\u00f73=\" }\n    @{ Row = 5; Col = 4; OldText = \"87\u00f74=\"; NewText = \"43\u00f79=\" }\n    @{ Row = 5; Col = 5; OldText = \"36\u00f75=\"; NewText = \"33\u00f73=\" }\n    @{ Row = 9; Col = 1; OldText = \"28\u00f74=\"; NewText = \"15\u00f75=\" }\n    @{ Row = 9; Col = 2; OldText = \"82\u00f76=\"; NewText = \"99\u00f72=\" }\n    @{ Row = 9; Col = 3; OldText = \"70\u00f74=\"; NewText = \"78\u00f72=\" }\n    @{ Row = 9; Col = 4; OldText = \"62\u00f73=\"; NewText = \"66\u00f78=\" }\n    @{ Row = 9; Col = 5; OldText = \"36\u00f77=\"; NewText = \"88\u00f74=\" }\n    @{ Row = 13; Col = 1; OldText = \"57\u00f75=\"; NewText = \"96\u00f72=\" }\n    @{ Row = 13; Col = 2; OldText = \"88\u00f74=\"; NewText = \"28\u00f78=\" }\n    @{ Row = 13; Col = 3; OldText = \"74\u00f72=\"; NewText = \"73\u00f77=\" }\n    @{ Row = 13; Col = 4; OldText = \"68\u00f77=\"; NewText = \"81\u00f73=\" }\n    @{ Row = 13; Col = 5; OldText = \"90\u00f72=\"; NewText = \"75\u00f73=\" }\n    @{ Row = 17; Col = 1; OldText = \"26\u00f76=\"; NewText = \"33\u00f77=\" }\n    @{ Row = 17; Col = 2; OldText = \"28\u00f76=\"; NewText = \"54\u00f79=\" }\n    @{ Row = 17; Col = 3; OldText = \"96\u00f73=\"; NewText = \"69\u00f76=\" }\n    @{ Row = 17; Col = 4; OldText = \"44\u00f79=\"; NewText = \"47\u00f76=\" }\n    @{ Row = 17; Col = 5; OldText = \"41\u00f74=\"; NewText = \"52\u00f76=\" }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $table.Cell($rep.Row, $rep.Col)\n    # Table-cell ranges report their text with a trailing \"end of cell\" mark\n    # (CR + cell-mark characters), so trim any trailing control characters\n    # before comparing against the plain expected text.\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -ne $rep.OldText) {\n        throw \"Unexpected cell text at Row=$($rep.Row) Col=$($rep.Col): expected '$($rep.OldText)' but found '$currentText'\"\n    }\n    $cell.Range.Text = $rep.NewText\n}"}
